$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: bugs re-tested and now safe ("Aman") ---
$ws.Range("D10").Value = "Aman"
$ws.Range("D16").Value = "Aman"
$ws.Range("D17").Value = "Aman"

# --- New role-testing blocks appended below the existing data (rows 46-74) ---
# Row 46
$ws.Range("A46").Value = "Manager"
$ws.Range("B46").Value = "Employee"
$ws.Range("C46").Value = "Add"
$ws.Range("D46").Value = "Format No Hp tidak valid ketika lebih dari 10 angka"
# Row 47
$ws.Range("C47").Value = "Edit"
$ws.Range("D47").Value = "Aman"
# Row 48
$ws.Range("C48").Value = "Delete"
$ws.Range("D48").Value = "Aman"
# Row 49
$ws.Range("C49").Value = "GetById"
$ws.Range("D49").Value = "Aman"
# Row 50
$ws.Range("C50").Value = "Get"
$ws.Range("D50").Value = "Error : The ObjectContext instance has been disposed …"
# Row 52
$ws.Range("B52").Value = "Region"
$ws.Range("C52").Value = "Add"
$ws.Range("D52").Value = "Aman"
# Row 53
$ws.Range("C53").Value = "Edit"
$ws.Range("D53").Value = "Aman"
# Row 54
$ws.Range("C54").Value = "Delete"
$ws.Range("D54").Value = "Aman"
# Row 55
$ws.Range("C55").Value = "GetById"
$ws.Range("D55").Value = "Aman"
# Row 56
$ws.Range("C56").Value = "Get"
$ws.Range("D56").Value = "Aman"
# Row 58
$ws.Range("B58").Value = "SubDistrict"
$ws.Range("C58").Value = "Add"
$ws.Range("D58").Value = "Aman"
# Row 59
$ws.Range("C59").Value = "Edit"
$ws.Range("D59").Value = "Aman"
# Row 60
$ws.Range("C60").Value = "Delete"
$ws.Range("D60").Value = "Aman"
# Row 61
$ws.Range("C61").Value = "GetById"
$ws.Range("D61").Value = "Aman"
# Row 62
$ws.Range("C62").Value = "Get"
$ws.Range("D62").Value = "Aman"
# Row 64
$ws.Range("A64").Value = "Admin Produksi"
$ws.Range("B64").Value = "Goods"
$ws.Range("C64").Value = "Add"
# Row 65
$ws.Range("C65").Value = "Edit"
# Row 66
$ws.Range("C66").Value = "Delete"
# Row 67
$ws.Range("C67").Value = "GetById"
# Row 68
$ws.Range("C68").Value = "Get"
$ws.Range("D68").Value = "Error : The ObjectContext instance has been disposed and can no longer be used for operations that require a connection."
# Row 70
$ws.Range("B70").Value = "Category"
$ws.Range("C70").Value = "Add"
$ws.Range("D70").Value = "Aman"
# Row 71
$ws.Range("C71").Value = "Edit"
$ws.Range("D71").Value = "Aman"
# Row 72
$ws.Range("C72").Value = "Delete"
$ws.Range("D72").Value = "Aman"
# Row 73
$ws.Range("C73").Value = "GetById"
$ws.Range("D73").Value = "Aman"
# Row 74
$ws.Range("C74").Value = "Get"
$ws.Range("D74").Value = "Aman"

# --- Column A widened to fit the new longer role name ("Admin Produksi") ---
$ws.Columns.Item(1).ColumnWidth = 16.9

# --- Reset the view: scroll back to top, select the last-edited block ---
$ws.Range("A1").Select() | Out-Null
$ws.Range("C70:D74").Select() | Out-Null
